# Automatic update: refreshed meteocat extraction timestamps + sensor
# readings (2026-02-05 16:49 run). Percentage-looking readings (column H)
# are written with a leading apostrophe so Excel keeps them as literal text
# ("70%") instead of auto-converting them to a numeric percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-05 16:47:42'
$ws.Range('K2').Value = '3.5 MJ/m2'
$ws.Range('E3').Value = '2026-02-05 16:47:44'
$ws.Range('K3').Value = '6.0 MJ/m2'
$ws.Range('O3').Value = '-2.0 °C'
$ws.Range('E4').Value = '2026-02-05 16:47:46'
$ws.Range('H4').Value = '''70%'
$ws.Range('J4').Value = '990.8 hPa'
$ws.Range('K4').Value = '6.6 MJ/m2'
$ws.Range('L4').Value = '79.6 km/h - 283º 16:27 TU'
$ws.Range('O4').Value = '10.7 °C'
$ws.Range('E5').Value = '2026-02-05 16:47:49'
$ws.Range('J5').Value = '990.8 hPa'
$ws.Range('K5').Value = '6.7 MJ/m2'
$ws.Range('O5').Value = '9.0 °C'
$ws.Range('E6').Value = '2026-02-05 16:47:51'
$ws.Range('J6').Value = '992.3 hPa'
$ws.Range('K6').Value = '3.7 MJ/m2'
$ws.Range('M6').Value = '15.6 °C 16:10 TU'
$ws.Range('O6').Value = '12.4 °C'
$ws.Range('E7').Value = '2026-02-05 16:47:54'
$ws.Range('J7').Value = '992.2 hPa'
$ws.Range('K7').Value = '3.0 MJ/m2'
$ws.Range('O7').Value = '10.1 °C'
$ws.Range('E8').Value = '2026-02-05 16:47:56'
$ws.Range('H8').Value = '''89%'
$ws.Range('K8').Value = '5.5 MJ/m2'
$ws.Range('O8').Value = '7.8 °C'
$ws.Range('E9').Value = '2026-02-05 16:47:58'
$ws.Range('O9').Value = '1.8 °C'
$ws.Range('E10').Value = '2026-02-05 16:48:01'
$ws.Range('H10').Value = '''91%'
$ws.Range('O10').Value = '6.9 °C'
$ws.Range('E11').Value = '2026-02-05 16:48:03'
$ws.Range('J11').Value = '995.7 hPa'
$ws.Range('K11').Value = '3.1 MJ/m2'
$ws.Range('O11').Value = '0.0 °C'
$ws.Range('E12').Value = '2026-02-05 16:48:06'
$ws.Range('K12').Value = '3.9 MJ/m2'
$ws.Range('M12').Value = '15.6 °C 16:04 TU'
$ws.Range('O12').Value = '9.0 °C'
$ws.Range('E13').Value = '2026-02-05 16:48:08'
$ws.Range('H13').Value = '''84%'
$ws.Range('O13').Value = '7.3 °C'
$ws.Range('E14').Value = '2026-02-05 16:48:10'
$ws.Range('H14').Value = '''69%'
$ws.Range('I14').Value = '5.2 mm'
$ws.Range('K14').Value = '1.7 MJ/m2'
$ws.Range('E15').Value = '2026-02-05 16:48:13'
$ws.Range('H15').Value = '''87%'
$ws.Range('J15').Value = '991.3 hPa'
$ws.Range('K15').Value = '6.6 MJ/m2'
$ws.Range('M15').Value = '16.1 °C 16:26 TU'
$ws.Range('O15').Value = '6.9 °C'
$ws.Range('E16').Value = '2026-02-05 16:48:15'
$ws.Range('K16').Value = '2.1 MJ/m2'
$ws.Range('L16').Value = '44.3 km/h - 230º 16:07 TU'
$ws.Range('M16').Value = '6.1 °C 16:18 TU'
$ws.Range('O16').Value = '3.4 °C'
$ws.Range('E17').Value = '2026-02-05 16:48:18'
$ws.Range('I17').Value = '8.0 mm'
$ws.Range('J17').Value = '995.6 hPa'
$ws.Range('E18').Value = '2026-02-05 16:48:20'
$ws.Range('K18').Value = '1.3 MJ/m2'
$ws.Range('E19').Value = '2026-02-05 16:48:22'
$ws.Range('I19').Value = '7.6 mm'
$ws.Range('J19').Value = '992.7 hPa'
$ws.Range('K19').Value = '4.8 MJ/m2'
$ws.Range('O19').Value = '7.1 °C'
$ws.Range('E20').Value = '2026-02-05 16:48:25'
$ws.Range('H20').Value = '''72%'
$ws.Range('K20').Value = '1.2 MJ/m2'
$ws.Range('E21').Value = '2026-02-05 16:48:27'
$ws.Range('H21').Value = '''85%'
$ws.Range('J21').Value = '991.7 hPa'
$ws.Range('K21').Value = '5.9 MJ/m2'
$ws.Range('M21').Value = '13.5 °C 16:11 TU'
$ws.Range('O21').Value = '5.4 °C'
$ws.Range('E22').Value = '2026-02-05 16:48:30'
$ws.Range('K22').Value = '5.2 MJ/m2'
$ws.Range('O22').Value = '7.6 °C'
$ws.Range('E23').Value = '2026-02-05 16:48:32'
$ws.Range('H23').Value = '''86%'
$ws.Range('J23').Value = '990.7 hPa'
$ws.Range('K23').Value = '3.4 MJ/m2'
$ws.Range('O23').Value = '7.8 °C'
$ws.Range('E24').Value = '2026-02-05 16:48:35'
$ws.Range('J24').Value = '989.8 hPa'
$ws.Range('K24').Value = '5.2 MJ/m2'
$ws.Range('O24').Value = '10.0 °C'
$ws.Range('E25').Value = '2026-02-05 16:48:37'
$ws.Range('I25').Value = '6.7 mm'
$ws.Range('J25').Value = '994.8 hPa'
$ws.Range('K25').Value = '3.2 MJ/m2'
$ws.Range('O25').Value = '0.3 °C'
$ws.Range('E26').Value = '2026-02-05 16:48:40'
$ws.Range('K26').Value = '4.2 MJ/m2'
$ws.Range('O26').Value = '-1.1 °C'
$ws.Range('E27').Value = '2026-02-05 16:48:42'
$ws.Range('J27').Value = '991.1 hPa'
$ws.Range('K27').Value = '3.5 MJ/m2'
$ws.Range('O27').Value = '8.0 °C'
$ws.Range('E28').Value = '2026-02-05 16:48:44'
$ws.Range('J28').Value = '994.0 hPa'
$ws.Range('O28').Value = '1.8 °C'
$ws.Range('E29').Value = '2026-02-05 16:48:47'
$ws.Range('H29').Value = '''85%'
$ws.Range('K29').Value = '4.5 MJ/m2'
$ws.Range('L29').Value = '68.8 km/h - 259º 16:28 TU'
$ws.Range('M29').Value = '15.8 °C 16:15 TU'
$ws.Range('O29').Value = '8.0 °C'
$ws.Range('E30').Value = '2026-02-05 16:48:49'
$ws.Range('H30').Value = '''67%'
$ws.Range('I30').Value = '4.9 mm'
$ws.Range('K30').Value = '1.3 MJ/m2'
$ws.Range('O30').Value = '-2.2 °C'
$ws.Range('E31').Value = '2026-02-05 16:48:51'
$ws.Range('I31').Value = '17.3 mm'
$ws.Range('J31').Value = '994.8 hPa'
$ws.Range('O31').Value = '3.7 °C'
$ws.Range('E32').Value = '2026-02-05 16:48:54'
$ws.Range('H32').Value = '''84%'
$ws.Range('J32').Value = '992.1 hPa'
$ws.Range('K32').Value = '6.1 MJ/m2'
$ws.Range('O32').Value = '11.5 °C'
$ws.Range('E33').Value = '2026-02-05 16:48:56'
$ws.Range('O33').Value = '8.2 °C'
$ws.Range('E34').Value = '2026-02-05 16:48:58'
$ws.Range('K34').Value = '2.6 MJ/m2'
$ws.Range('L34').Value = '34.2 km/h - 254º 16:06 TU'
$ws.Range('M34').Value = '10.1 °C 16:12 TU'
$ws.Range('O34').Value = '2.8 °C'
$ws.Range('E35').Value = '2026-02-05 16:49:01'
$ws.Range('I35').Value = '3.4 mm'
$ws.Range('K35').Value = '2.8 MJ/m2'
$ws.Range('M35').Value = '-1.2 °C 16:07 TU'
$ws.Range('O35').Value = '-3.1 °C'
$ws.Range('E36').Value = '2026-02-05 16:49:03'
$ws.Range('J36').Value = '992.6 hPa'
$ws.Range('K36').Value = '9.6 MJ/m2'
$ws.Range('O36').Value = '9.7 °C'
